# data : case 1
# Update the values in column A (A1:A33) of the active worksheet to the
# new values described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    0.99999999487151059,
    0.99889873536450113,
    0.99656729392587118,
    1.0016924836169485,
    0.99375030572468159,
    0.97569977940575847,
    0.97210189935938152,
    0.96705186783595831,
    0.96252199504755964,
    0.95891390538645804,
    0.95842214483027766,
    0.95564312960357634,
    0.94435875307966199,
    0.94019194994113042,
    0.93760072266518246,
    0.93509441104433577,
    0.93138664496422197,
    0.93027773140235026,
    0.99784551130628985,
    0.99072844714055697,
    0.98932994630757998,
    0.98806543754887777,
    0.98823368544268675,
    0.97521354860285159,
    0.96875667218342998,
    0.97403454367679687,
    0.97196329154095462,
    0.96448720050794035,
    0.9596405234292511,
    0.958302634499689,
    0.96394107854377409,
    0.96688269697347629,
    0.96636268233891209
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
